$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''57.262.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.80%  '

$ws.Range('D3').Value = '''3.262.84'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.59%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '''397.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.33%  '

$ws.Range('D6').Value = '''108.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.75%  '

$ws.Range('D7').Value = '''0.581'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.94%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').Value = '''0.621'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.07%  '

$ws.Range('D10').Value = '''39.45'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.57%  '

$ws.Range('D11').Value = '''0.0956'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.24%  '

$ws.Range('E12').Value = '  +1.68%  '

$ws.Range('D13').Value = '''3.781.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.73%  '

$ws.Range('E14').Value = '  +2.89%  '

$ws.Range('E15').Value = '  -0.09%  '

$ws.Range('D16').Value = '''3.383.95'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.05%  '

$ws.Range('E17').Value = '  -1.47%  '

$ws.Range('D18').Value = '''11.09'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.70%  '

$ws.Range('D19').Value = '''57.051.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.60%  '

$ws.Range('E20').Value = '  -1.45%  '

$ws.Range('E21').Value = '  +7.11%  '

$ws.Range('D22').Value = '''12.96'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.09%  '

$ws.Range('D23').Value = '''294.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.88%  '

$ws.Range('D24').Value = '''74.24'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.56%  '

$ws.Range('D25').Value = '''3.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.50%  '

$ws.Range('B26').Value = 'Filecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D26').Value = '''7.98'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.54%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''28.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.05%  '

$ws.Range('E28').Value = '  +0.50%  '

$ws.Range('D29').Value = '''7.41'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E30').Value = '  -2.84%  '

$ws.Range('E32').Value = '  +1.06%  '

$ws.Range('D33').Value = '''11.22'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.72%  '

$ws.Range('D34').Value = '''40.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +11.07%  '

$ws.Range('D35').Value = '''0.0487'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.25%  '

$ws.Range('D36').Value = '''2.14'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.23%  '

$ws.Range('D37').Value = '''51.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.30%  '

$ws.Range('E38').Value = '  -0.14%  '

$ws.Range('E39').Value = '  -0.78%  '

$ws.Range('D40').Value = '''3.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.29%  '

$ws.Range('D41').Value = '''136.70'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.63%  '

$ws.Range('E42').Value = '  +1.56%  '

$ws.Range('E43').Value = '  -2.50%  '

$ws.Range('E44').Value = '  -2.51%  '

$ws.Range('E45').Value = '  -1.24%  '

$ws.Range('D46').Value = '''16.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.42%  '

$ws.Range('D47').Value = '''22.36'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.78%  '

$ws.Range('E48').Value = '  +4.89%  '

$ws.Range('D49').Value = '''2.150.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.43%  '

$ws.Range('D50').Value = '''1.98'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.14%  '

$ws.Range('E51').Value = '  -6.17%  '
